$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stored username (A2) and password (B2).
$ws.Range("A2").Value = "venukollapudi@gmail.com"
$ws.Range("B2").Value = "Venu@12345"

# The new password value contains an "@", so Excel auto-recognizes it as an
# email address and turns it into a mailto hyperlink (as already happened
# for the username in A2), applying the built-in "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Venu@12345")
$ws.Range("B2").Style = "Hyperlink"

# The active cell / selection moved to A2.
$ws.Range("A2").Select()
